# Apply the "Office Theme" design colours to the presentation's slide
# master theme (ppt/theme/theme1.xml), and switch every table in the deck
# from the old custom "Table_0" style to the built-in Office table style.

$p = $ppt.ActivePresentation

# --- 1. Re-colour the design theme (was "Integral" / Red Violet) to the
#        stock Office Theme colour scheme. -------------------------------
$master = $p.SlideMaster
$theme  = $master.Theme
$tcs    = $theme.ThemeColorScheme

# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (VBA RGB() little-endian
# BGR packing: R + G*256 + B*65536) for the standard "Office" colour scheme.
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}

# --- 2. Point every table at the built-in table style instead of the
#        custom "Table_0" style that used to live in tableStyles.xml. ----
$newTableStyle = "{207AA8BA-5E77-4AA8-85F0-E5E3DC3E65DB}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTable) {
            $sh.Table.ApplyStyle($newTableStyle)
        }
    }
}
